$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to be bumped
# from 46060 (2026-02-07) to 46061 (2026-02-08) for every data row (2..103).
for ($row = 2; $row -le 103; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
